$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 69, shifting existing rows 69:139 down to 70:140
$ws.Rows.Item(69).Insert()

# Populate the new row 69 with the new weekly record.
# Columns that stay constant (copied from the pattern used by the rest of the dataset):
$ws.Range("A69").Value = 5
$ws.Range("B69").Value = "Macroferia Regional de Talca"
$ws.Range("C69").Value = "Maule"
$ws.Range("E69").Value = 7
$ws.Range("F69").Value = 300000000
$ws.Range("G69").Value = "Espárragos"
$ws.Range("I69").Value = "Primera"
$ws.Range("N69").Value = "`$/kilo"
$ws.Range("Q69").Value = 1
$ws.Range("R69").Value = "Hortaliza"

# New values for the inserted record
$ws.Range("D69").Value = 45280
$ws.Range("H69").Value = "Sin especificar"
$ws.Range("J69").Value = 3000
$ws.Range("K69").Value = 1300
$ws.Range("L69").Value = 1300
$ws.Range("M69").Value = 1300
$ws.Range("O69").Value = "Región del Maule"
$ws.Range("P69").Value = 1300
